$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, whether the text must be forced
# to remain plain text (needed when the string also parses as a number,
# which Excel would otherwise auto-convert).
$updates = @(
    @{ Cell = "D2"; Value = "30.504.70"; ForceText = 0 }
    @{ Cell = "E2"; Value = "  -0.98%  "; ForceText = 0 }
    @{ Cell = "D3"; Value = "1.911.73"; ForceText = 0 }
    @{ Cell = "E3"; Value = "  -1.58%  "; ForceText = 0 }
    @{ Cell = "D4"; Value = "0.9999"; ForceText = 1 }
    @{ Cell = "E4"; Value = "  -0.05%  "; ForceText = 0 }
    @{ Cell = "D5"; Value = "239.44"; ForceText = 1 }
    @{ Cell = "E5"; Value = "  -1.17%  "; ForceText = 0 }
    @{ Cell = "E6"; Value = "  +0.04%  "; ForceText = 0 }
    @{ Cell = "D7"; Value = "0.4777"; ForceText = 1 }
    @{ Cell = "E7"; Value = "  -2.36%  "; ForceText = 0 }
    @{ Cell = "D8"; Value = "0.2840"; ForceText = 1 }
    @{ Cell = "E8"; Value = "  -3.61%  "; ForceText = 0 }
    @{ Cell = "D9"; Value = "0.06710"; ForceText = 1 }
    @{ Cell = "E9"; Value = "  -2.70%  "; ForceText = 0 }
    @{ Cell = "D10"; Value = "18.79"; ForceText = 1 }
    @{ Cell = "E10"; Value = "  -3.30%  "; ForceText = 0 }
    @{ Cell = "D11"; Value = "101.88"; ForceText = 1 }
    @{ Cell = "E11"; Value = "  -4.26%  "; ForceText = 0 }
    @{ Cell = "D12"; Value = "1.915.02"; ForceText = 0 }
    @{ Cell = "E12"; Value = "  -1.37%  "; ForceText = 0 }
    @{ Cell = "D13"; Value = "0.07686"; ForceText = 1 }
    @{ Cell = "E13"; Value = "  -0.53%  "; ForceText = 0 }
    @{ Cell = "D14"; Value = "5.195"; ForceText = 1 }
    @{ Cell = "E14"; Value = "  -2.77%  "; ForceText = 0 }
    @{ Cell = "D15"; Value = "0.6714"; ForceText = 1 }
    @{ Cell = "E15"; Value = "  -4.04%  "; ForceText = 0 }
    @{ Cell = "D16"; Value = "30.514.30"; ForceText = 0 }
    @{ Cell = "E16"; Value = "  -0.93%  "; ForceText = 0 }
    @{ Cell = "D17"; Value = "259.66"; ForceText = 1 }
    @{ Cell = "E17"; Value = "  -6.26%  "; ForceText = 0 }
    @{ Cell = "D18"; Value = "1.000"; ForceText = 1 }
    @{ Cell = "E18"; Value = "  +0.00%  "; ForceText = 0 }
    @{ Cell = "D19"; Value = "0.000007475"; ForceText = 1 }
    @{ Cell = "E19"; Value = "  -3.16%  "; ForceText = 0 }
    @{ Cell = "D20"; Value = "12.66"; ForceText = 1 }
    @{ Cell = "E20"; Value = "  -3.45%  "; ForceText = 0 }
    @{ Cell = "D21"; Value = "5.387"; ForceText = 1 }
    @{ Cell = "E21"; Value = "  -0.99%  "; ForceText = 0 }
    @{ Cell = "D22"; Value = "1.001"; ForceText = 1 }
    @{ Cell = "E22"; Value = "  +0.00%  "; ForceText = 0 }
    @{ Cell = "D23"; Value = "6.278"; ForceText = 1 }
    @{ Cell = "E23"; Value = "  -3.71%  "; ForceText = 0 }
    @{ Cell = "D24"; Value = "9.359"; ForceText = 1 }
    @{ Cell = "E24"; Value = "  -3.58%  "; ForceText = 0 }
    @{ Cell = "D25"; Value = "167.66"; ForceText = 1 }
    @{ Cell = "E25"; Value = "  -0.32%  "; ForceText = 0 }
    @{ Cell = "D26"; Value = "19.15"; ForceText = 1 }
    @{ Cell = "E26"; Value = "  -2.58%  "; ForceText = 0 }
    @{ Cell = "D27"; Value = "2.062"; ForceText = 1 }
    @{ Cell = "E27"; Value = "  -4.32%  "; ForceText = 0 }
    @{ Cell = "D28"; Value = "4.814"; ForceText = 1 }
    @{ Cell = "E28"; Value = "  +5.81%  "; ForceText = 0 }
    @{ Cell = "D29"; Value = "1.382"; ForceText = 1 }
    @{ Cell = "E29"; Value = "  -0.38%  "; ForceText = 0 }
    @{ Cell = "D30"; Value = "0.09999"; ForceText = 1 }
    @{ Cell = "E30"; Value = "  -4.05%  "; ForceText = 0 }
    @{ Cell = "D32"; Value = "4.259"; ForceText = 1 }
    @{ Cell = "E32"; Value = "  -2.36%  "; ForceText = 0 }
    @{ Cell = "D33"; Value = "0.04724"; ForceText = 1 }
    @{ Cell = "E33"; Value = "  -2.67%  "; ForceText = 0 }
    @{ Cell = "D34"; Value = "0.7269"; ForceText = 1 }
    @{ Cell = "E34"; Value = "  -3.08%  "; ForceText = 0 }
    @{ Cell = "D35"; Value = "1.107"; ForceText = 1 }
    @{ Cell = "D36"; Value = "2.714"; ForceText = 1 }
    @{ Cell = "E36"; Value = "  -0.37%  "; ForceText = 0 }
    @{ Cell = "E37"; Value = "  -3.48%  "; ForceText = 0 }
    @{ Cell = "D38"; Value = "2.627"; ForceText = 1 }
    @{ Cell = "E38"; Value = "  -1.27%  "; ForceText = 0 }
    @{ Cell = "D39"; Value = "6.261"; ForceText = 1 }
    @{ Cell = "E39"; Value = "  -2.84%  "; ForceText = 0 }
    @{ Cell = "D40"; Value = "75.03"; ForceText = 1 }
    @{ Cell = "E40"; Value = "  -4.51%  "; ForceText = 0 }
    @{ Cell = "D41"; Value = "1.971"; ForceText = 1 }
    @{ Cell = "E41"; Value = "  -5.86%  "; ForceText = 0 }
    @{ Cell = "D42"; Value = "0.8615"; ForceText = 1 }
    @{ Cell = "E42"; Value = "  -4.98%  "; ForceText = 0 }
    @{ Cell = "D43"; Value = "105.62"; ForceText = 1 }
    @{ Cell = "E43"; Value = "  -2.80%  "; ForceText = 0 }
    @{ Cell = "D44"; Value = "0.4251"; ForceText = 1 }
    @{ Cell = "E44"; Value = "  -3.43%  "; ForceText = 0 }
    @{ Cell = "D45"; Value = "1.000"; ForceText = 1 }
    @{ Cell = "E45"; Value = "  +0.25%  "; ForceText = 0 }
    @{ Cell = "D46"; Value = "7.390"; ForceText = 1 }
    @{ Cell = "E46"; Value = "  -4.43%  "; ForceText = 0 }
    @{ Cell = "D47"; Value = "0.1202"; ForceText = 1 }
    @{ Cell = "E47"; Value = "  -3.33%  "; ForceText = 0 }
    @{ Cell = "D48"; Value = "919.48"; ForceText = 1 }
    @{ Cell = "E48"; Value = "  -6.58%  "; ForceText = 0 }
    @{ Cell = "D49"; Value = "34.74"; ForceText = 1 }
    @{ Cell = "E49"; Value = "  -3.23%  "; ForceText = 0 }
    @{ Cell = "D50"; Value = "8.814"; ForceText = 1 }
    @{ Cell = "E50"; Value = "  -4.97%  "; ForceText = 0 }
    @{ Cell = "D51"; Value = "0.05745"; ForceText = 1 }
    @{ Cell = "E51"; Value = "  +0.25%  "; ForceText = 0 }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    if ($u.ForceText -eq 1) {
        $c.NumberFormat = "@"
        $c.Value = $u.Value
        $c.Style = "Normal"
    } else {
        $c.Value = $u.Value
    }
}
